$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Underlying data edits (DM2/ia.xlsx, Sheet1):
#   D4 (trip/Ca)          : 1 -> 2
#   E5 (usersessions/Na)  : 4 -> 1
#   D6 (utils/Ca)         : 1 -> 3
$ws.Range("D4").Value = 2
$ws.Range("E5").Value = 1
$ws.Range("D6").Value = 3

# Re-enter the B/C formulas across the whole B4:B10 / C4:C10 block so Excel
# regroups them into shared formulas (as it does when a formula is typed
# once and propagated down a contiguous range), matching the refreshed
# "si" shared-formula groups seen after the edit.
$ws.Range("B4:B10").Formula = "=D4/(E4+D4)"
$ws.Range("C4:C10").Formula = "=F4/G4"

# Recalculate so cached formula results reflect the new inputs.
$excel.CalculateFull()

# Selection moved to E4 after the edit.
$ws.Range("E4").Select()
